$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate "Период" / [date]->[sum_amount] column (column F),
# which shifts the old column G ("Конечный остаток" / [saldo_end]) left into F.
$ws.Range("F1").EntireColumn.Delete()

# Rename the date placeholder to [title]
$ws.Range("B6").Value = "[title]"

# Match the author's final selection/cursor position
[void]$ws.Range("H8").Select()
